$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the pickup/dropoff date values (drives the sharedStrings churn) ---
$ws.Range("C2").Value2 = "12/10/2017"
$ws.Range("D2").Value2 = "12/24/2017"

# --- Move the active selection to D2 (matches sheetView/selection in the diff) ---
$ws.Range("D2").Select()
